$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New fictitious surnames (column A) and first names (column C) + sex (column U)
# for rows 2-18. Column B only has a value on row 2 ("ANDRE"), already correct.
$data = @(
    @{ Row = 2;  Nom = "PERRIN";   Prenom = "François";  Sexe = "M" },
    @{ Row = 3;  Nom = "ROBIN";    Prenom = "Valérie";    Sexe = "F" },
    @{ Row = 4;  Nom = "CLEMENT";  Prenom = "Sandrine";   Sexe = "F" },
    @{ Row = 5;  Nom = "MORIN";    Prenom = "Stéphanie";  Sexe = "F" },
    @{ Row = 6;  Nom = "NICOLAS";  Prenom = "Sophie";     Sexe = "F" },
    @{ Row = 7;  Nom = "HENRY";    Prenom = "Véronique";  Sexe = "F" },
    @{ Row = 8;  Nom = "ROUSSEL";  Prenom = "Céline";     Sexe = "F" },
    @{ Row = 9;  Nom = "MATHIEU";  Prenom = "Chantal";    Sexe = "F" },
    @{ Row = 10; Nom = "GAUTIER";  Prenom = "Christiane"; Sexe = "F" },
    @{ Row = 11; Nom = "MASSON";   Prenom = "Didier";     Sexe = "M" },
    @{ Row = 12; Nom = "MARCHAND"; Prenom = "René";       Sexe = "M" },
    @{ Row = 13; Nom = "DUVAL";    Prenom = "Vincent";    Sexe = "M" },
    @{ Row = 14; Nom = "DENIS";    Prenom = "Jeanne";     Sexe = "F" },
    @{ Row = 15; Nom = "DUMONT";   Prenom = "Patricia";   Sexe = "F" },
    @{ Row = 16; Nom = "MARIE";    Prenom = "Guillaume";  Sexe = "M" },
    @{ Row = 17; Nom = "LEMAIRE";  Prenom = "Annie";      Sexe = "F" },
    @{ Row = 18; Nom = "NOEL";     Prenom = "Bruno";      Sexe = "M" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.Nom
    $ws.Range("C$r").Value = $item.Prenom
    $ws.Range("U$r").Value = $item.Sexe
}

# B2 keeps its existing text ("ANDRE"); B3:B18 stay blank.

# Re-apply the formatting (style) used by A2:C18 and U2:U18 so it matches the
# style already used elsewhere on the sheet (numeric format 14 / centered).
$ws.Range("V2").Copy()
$ws.Range("A2:C18").PasteSpecial(-4122)
$ws.Range("U2:U18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the remembered selection on the sheet.
$ws.Range("S13").Select()
